$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Doppia iscrizione"

$ws.Range("B10").Select()
